$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.782.89"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "3.558.76"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "653.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "3.555.37"
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.22%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "4.220.58"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").Value = "94.779.08"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "3.560.20"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.27%  "
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "506.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("E24").Value = "  -4.78%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -2.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "94.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.03%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "3.748.76"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.29%  "
$ws.Range("E31").Value = "  -1.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("E35").Value = "  -4.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.49%  "
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "579.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.901"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.70%  "
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "34.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +31.26%  "
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0412"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -1.54%  "
